$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fix the "Periodo Mora" (E column) ordering and the matching "Valor
#    Mora" (F column) values on rows 16-24 - the source data had a few
#    periods out of order / their amounts swapped.
# ---------------------------------------------------------------------------
$ws.Range("E16").Value = "1708"
$ws.Range("E17").Value = "1709"
$ws.Range("E18").Value = "1711"
$ws.Range("F18").Value = 11733
$ws.Range("E19").Value = "1712"
$ws.Range("F19").Value = 32000
$ws.Range("E20").Value = "1801"
$ws.Range("F20").Value = 27600
$ws.Range("E21").Value = "1802"
$ws.Range("E23").Value = "1804"
$ws.Range("E24").Value = "1805"
$ws.Range("F24").Value = 36000

# ---------------------------------------------------------------------------
# 2) Insert 10 fresh rows right after the current last data row (25) so the
#    footer ("NOMBRE DEL REPRESENTANTE LEGAL" / "FIRMA DEL REPRESENTANTE
#    LEGAL" block, previously rows 30:31) is pushed down to rows 40:41.
# ---------------------------------------------------------------------------
$ws.Rows("25:34").Insert()

# ---------------------------------------------------------------------------
# 3) Populate the newly inserted rows 26:34 with a copy (values + format) of
#    the corrected data rows 16:24 - this is "parte 1 de nuevos estado de
#    cuenta" being appended to the worker table.
# ---------------------------------------------------------------------------
$ws.Range("B16:J24").Copy($ws.Range("B26:J34"))

# ---------------------------------------------------------------------------
# 4) The row that used to be row 25 (special bold/shaded "last row" style)
#    got pushed down to row 35 by the insert, carrying its own formatting
#    along with it. Recreate its data - with the normal row 24 formatting -
#    in the now-empty row 25 so that it reads as a regular body row, exactly
#    like the other entries before the true, still-special, final row (35).
# ---------------------------------------------------------------------------
$ws.Range("B24:J24").Copy()
$ws.Range("B25:J25").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("B25").Value = "CC"
$ws.Range("C25").Value = "1002241811"
$ws.Range("D25").Value = "EUDENIS PATRICIA OLIVO JIMENEZ"
$ws.Range("E25").Value = "2007"
$ws.Range("F25").Value = 25200
$ws.Range("G25").Value = 900000

# ---------------------------------------------------------------------------
# 5) Update the "Valor Mora" grand total shown near the top of the sheet -
#    the table now holds twice as many rows of mora values.
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 609066

$excel.CutCopyMode = 0
